# Applies the commit: "Added new files for equity calculation"
# 1. Re-format the existing "trade_date" column (I2:I361) to the
#    YYYY-MM-DD HH:MM:SS datetime format (matching column B's style),
#    while the four freshly appended rows keep the original
#    YYYY-MM-DD date-only format for that column.
# 2. Append 4 new trading-day rows (362-365) of Suzlon Energy Ltd
#    NSE history data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: re-style the existing trade_date column (I2:I361) ---
$ws.Range("I2:I361").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Step 2: append the 4 new data rows ---
$ws.Range("A362").Value = 45.79
$ws.Range("B362").Value = 46049
$ws.Range("B362").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C362").Value = "NSE"
$ws.Range("D362").Value = 46.25
$ws.Range("E362").Value = 44.88
$ws.Range("F362").Value = 45.84
$ws.Range("G362").Value = "SUZENE"
$ws.Range("H362").Value = 62382536
$ws.Range("I362").Value = 46049
$ws.Range("I362").NumberFormat = "YYYY-MM-DD"
$ws.Range("J362").Value = "INE040H01021"
$ws.Range("K362").Value = "Suzlon Engergy Ltd"
$ws.Range("L362").Value = "SUZENE"
$ws.Range("M362").Value = "BREEZE"

$ws.Range("A363").Value = 47.8
$ws.Range("B363").Value = 46050
$ws.Range("B363").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C363").Value = "NSE"
$ws.Range("D363").Value = 48
$ws.Range("E363").Value = 46.15
$ws.Range("F363").Value = 46.15
$ws.Range("G363").Value = "SUZENE"
$ws.Range("H363").Value = 59757766
$ws.Range("I363").Value = 46050
$ws.Range("I363").NumberFormat = "YYYY-MM-DD"
$ws.Range("J363").Value = "INE040H01021"
$ws.Range("K363").Value = "Suzlon Engergy Ltd"
$ws.Range("L363").Value = "SUZENE"
$ws.Range("M363").Value = "BREEZE"

$ws.Range("A364").Value = 47.44
$ws.Range("B364").Value = 46051
$ws.Range("B364").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C364").Value = "NSE"
$ws.Range("D364").Value = 48.03
$ws.Range("E364").Value = 46.63
$ws.Range("F364").Value = 48.03
$ws.Range("G364").Value = "SUZENE"
$ws.Range("H364").Value = 54843475
$ws.Range("I364").Value = 46051
$ws.Range("I364").NumberFormat = "YYYY-MM-DD"
$ws.Range("J364").Value = "INE040H01021"
$ws.Range("K364").Value = "Suzlon Engergy Ltd"
$ws.Range("L364").Value = "SUZENE"
$ws.Range("M364").Value = "BREEZE"

$ws.Range("A365").Value = 47.67
$ws.Range("B365").Value = 46052
$ws.Range("B365").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C365").Value = "NSE"
$ws.Range("D365").Value = 47.95
$ws.Range("E365").Value = 46.81
$ws.Range("F365").Value = 47.4
$ws.Range("G365").Value = "SUZENE"
$ws.Range("H365").Value = 60470060
$ws.Range("I365").Value = 46052
$ws.Range("I365").NumberFormat = "YYYY-MM-DD"
$ws.Range("J365").Value = "INE040H01021"
$ws.Range("K365").Value = "Suzlon Engergy Ltd"
$ws.Range("L365").Value = "SUZENE"
$ws.Range("M365").Value = "BREEZE"

Write-Output "Applied equity history update: re-styled I2:I361 and appended rows 362-365"
